$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (e.g. "1.011", "28.672.34") that must
# stay literal text rather than being auto-converted to a Double by Excel.
# Temporarily force text format, assign, then strip the format again so the
# cell ends up with no explicit style (matching the original unstyled cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "28.672.34"
$ws.Range("E2").Value = "  -3.03%  "
Set-TextValue $ws.Range("D3") "1.955.59"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "321.59"
$ws.Range("E5").Value = "  -2.47%  "
Set-TextValue $ws.Range("D6") "1.011"
$ws.Range("E6").Value = "  +0.05%  "
Set-TextValue $ws.Range("D7") "0.4774"
$ws.Range("E7").Value = "  -4.57%  "
Set-TextValue $ws.Range("D8") "0.4051"
$ws.Range("E8").Value = "  -3.95%  "
Set-TextValue $ws.Range("D9") "53.81"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -6.09%  "
Set-TextValue $ws.Range("D11") "1.063"
$ws.Range("E11").Value = "  -4.91%  "
Set-TextValue $ws.Range("D12") "22.30"
$ws.Range("E12").Value = "  -4.51%  "
Set-TextValue $ws.Range("D13") "1.924.52"
$ws.Range("E13").Value = "  -4.53%  "
Set-TextValue $ws.Range("D14") "7.636"
$ws.Range("E14").Value = "  -5.36%  "
Set-TextValue $ws.Range("D15") "6.227"
$ws.Range("E15").Value = "  -3.87%  "
Set-TextValue $ws.Range("D16") "1.014"
$ws.Range("E16").Value = "  +0.14%  "
Set-TextValue $ws.Range("D17") "0.00001077"
$ws.Range("E17").Value = "  -3.15%  "
Set-TextValue $ws.Range("D18") "89.35"
$ws.Range("E18").Value = "  -5.03%  "
Set-TextValue $ws.Range("D19") "0.06634"
$ws.Range("E19").Value = "  -0.56%  "
Set-TextValue $ws.Range("D20") "18.71"
$ws.Range("E20").Value = "  -5.07%  "
Set-TextValue $ws.Range("D21") "1.012"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -2.06%  "
Set-TextValue $ws.Range("D23") "28.663.55"
$ws.Range("E23").Value = "  -3.04%  "
Set-TextValue $ws.Range("D24") "11.56"
$ws.Range("E24").Value = "  -3.40%  "
Set-TextValue $ws.Range("D25") "2.294"
$ws.Range("E25").Value = "  -0.47%  "
Set-TextValue $ws.Range("D26") "2.170.18"
$ws.Range("E26").Value = "  -3.77%  "
Set-TextValue $ws.Range("D27") "154.51"
$ws.Range("E27").Value = "  -2.53%  "
Set-TextValue $ws.Range("D28") "20.23"
$ws.Range("E28").Value = "  -2.16%  "
Set-TextValue $ws.Range("D29") "5.961"
$ws.Range("E29").Value = "  -7.16%  "
Set-TextValue $ws.Range("D30") "2.159"
$ws.Range("E30").Value = "  -6.14%  "
Set-TextValue $ws.Range("D31") "124.08"
$ws.Range("E31").Value = "  -3.25%  "
Set-TextValue $ws.Range("D32") "1.004"
$ws.Range("E32").Value = "  -4.57%  "
Set-TextValue $ws.Range("D33") "0.09593"
$ws.Range("E33").Value = "  -3.41%  "
Set-TextValue $ws.Range("D34") "5.657"
$ws.Range("E34").Value = "  -2.95%  "
Set-TextValue $ws.Range("D35") "1.442"
$ws.Range("E35").Value = "  -8.38%  "
Set-TextValue $ws.Range("D36") "3.668"
$ws.Range("E36").Value = "  -3.39%  "
Set-TextValue $ws.Range("D37") "0.02349"
$ws.Range("E37").Value = "  -4.87%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "1.272"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.06239"
$ws.Range("E39").Value = "  -1.72%  "
Set-TextValue $ws.Range("D40") "8.732"
$ws.Range("E40").Value = "  -6.03%  "
Set-TextValue $ws.Range("D41") "0.6243"
$ws.Range("E41").Value = "  -4.91%  "
Set-TextValue $ws.Range("D42") "11.12"
$ws.Range("E42").Value = "  -4.83%  "
Set-TextValue $ws.Range("D43") "1.011"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  -6.30%  "
Set-TextValue $ws.Range("D45") "1.334"
$ws.Range("E45").Value = "  +2.06%  "
Set-TextValue $ws.Range("D46") "0.5950"
$ws.Range("E46").Value = "  -6.23%  "
Set-TextValue $ws.Range("D47") "12.94"
$ws.Range("E47").Value = "  -3.82%  "
Set-TextValue $ws.Range("D48") "2.084"
$ws.Range("E48").Value = "  -5.27%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D49") "3.405"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D50") "0.00000000335"
$ws.Range("E50").Value = "  -2.80%  "
Set-TextValue $ws.Range("D51") "0.06830"
$ws.Range("E51").Value = "  -2.10%  "
